# Delete the "月はどこへ消えた？" post row (row 209).
# All rows below it shift up by one, which matches the target diff
# (the rest of the sheet is unchanged content, just renumbered).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(209).Delete()
